# Auto-update Data Telemetría - Ejecución Diaria $(date)
# Appends the new daily snapshot (fecha serial 45995) to the "historico"
# sheet and refreshes the "ultimo_snapshot" sheet with the same latest rows.

$wb = $excel.ActiveWorkbook

$historico = $wb.Worksheets.Item("historico")
$snapshot  = $wb.Worksheets.Item("ultimo_snapshot")

# New daily rows: [fecha, resumen, total_vin, cnt_Conectado 0-2,
#   cnt_Intermitente 3-14, cnt_Limitado 15-30+, cnt_Desconectado 31+,
#   cnt_Nunca, pct_Conectado 0-2, pct_Intermitente 3-14,
#   pct_Limitado 15-30+, pct_Desconectado 31+, pct_Nunca]
$newRows = @(
    @(45995, "Telemetría", 5904, 3545, 511, 187, 662, 999, 60.04, 8.66, 3.17, 11.21, 16.92),
    @(45995, "GPS (según REGLA)", 5300, 4667, 352, 102, 172, 7, 88.06, 6.64, 1.92, 3.25, 0.13),
    @(45995, "GPS (todas con gps_timestamp)", 11197, 9495, 835, 299, 568, 0, 84.8, 7.46, 2.67, 5.07, 0)
)

# Append the new rows to "historico", right after the last used row.
$usedRange = $historico.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1
$startRow = $lastRow + 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $data = $newRows[$i]
    for ($c = 0; $c -lt $data.Count; $c++) {
        $historico.Cells.Item($row, $c + 1).Value = $data[$c]
    }
    # Match the "fecha" column date formatting used by the existing rows.
    $historico.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD"
}

# Refresh "ultimo_snapshot" (rows 2-4) with the same latest rows.
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = 2 + $i
    $data = $newRows[$i]
    for ($c = 0; $c -lt $data.Count; $c++) {
        $snapshot.Cells.Item($row, $c + 1).Value = $data[$c]
    }
}

Write-Host "Applied daily update: historico rows $startRow-$($startRow + $newRows.Count - 1), ultimo_snapshot rows 2-4"
